$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 16 with the new "Applying feedback" entry
$ws.Range("A16").Value = "Applying feedback"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "2023-06-08"
$ws.Range("D16").Value = "Fixed search endpoint to apply none,1 or more filters. Returning the facility when creating and updating it.Removed some unnecessary comments, fixed pagination, added sanitization to all other endpoints, added model classes."

$ws.Rows.Item(16).RowHeight = 44.25

$excel.CalculateFull()

# Update sheet view: scroll position and selection
$ws.Range("J22").Select()
$excel.ActiveWindow.ScrollRow = 10
